$d = $word.ActiveDocument

$d.Content.Find.Execute("21+74=", $true, $false, $false, $false, $false, $true, 1, $false, "96-57=", 2) | Out-Null
$d.Content.Find.Execute("53+25=", $true, $false, $false, $false, $false, $true, 1, $false, "84+11=", 2) | Out-Null
$d.Content.Find.Execute("22-14=", $true, $false, $false, $false, $false, $true, 1, $false, "33-30=", 2) | Out-Null
$d.Content.Find.Execute("74-31=", $true, $false, $false, $false, $false, $true, 1, $false, "8+22=", 2) | Out-Null
$d.Content.Find.Execute("77-40=", $true, $false, $false, $false, $false, $true, 1, $false, "51-28=", 2) | Out-Null
$d.Content.Find.Execute("5+42=", $true, $false, $false, $false, $false, $true, 1, $false, "82-50=", 2) | Out-Null
$d.Content.Find.Execute("84-67=", $true, $false, $false, $false, $false, $true, 1, $false, "95-10=", 2) | Out-Null
$d.Content.Find.Execute("61-56=", $true, $false, $false, $false, $false, $true, 1, $false, "80-37=", 2) | Out-Null
$d.Content.Find.Execute("10-5=", $true, $false, $false, $false, $false, $true, 1, $false, "40-35=", 2) | Out-Null
$d.Content.Find.Execute("42+14=", $true, $false, $false, $false, $false, $true, 1, $false, "38+27=", 2) | Out-Null
$d.Content.Find.Execute("8+14=", $true, $false, $false, $false, $false, $true, 1, $false, "23+16=", 2) | Out-Null
$d.Content.Find.Execute("97-74=", $true, $false, $false, $false, $false, $true, 1, $false, "75-57=", 2) | Out-Null
$d.Content.Find.Execute("27-27=", $true, $false, $false, $false, $false, $true, 1, $false, "21+11=", 2) | Out-Null
$d.Content.Find.Execute("14+30=", $true, $false, $false, $false, $false, $true, 1, $false, "0+4=", 2) | Out-Null
$d.Content.Find.Execute("10+11=", $true, $false, $false, $false, $false, $true, 1, $false, "45+40=", 2) | Out-Null
$d.Content.Find.Execute("42+44=", $true, $false, $false, $false, $false, $true, 1, $false, "53-50=", 2) | Out-Null
$d.Content.Find.Execute("44-29=", $true, $false, $false, $false, $false, $true, 1, $false, "46-33=", 2) | Out-Null
$d.Content.Find.Execute("35-31=", $true, $false, $false, $false, $false, $true, 1, $false, "88-26=", 2) | Out-Null
$d.Content.Find.Execute("4+86=", $true, $false, $false, $false, $false, $true, 1, $false, "8+60=", 2) | Out-Null
$d.Content.Find.Execute("12+47=", $true, $false, $false, $false, $false, $true, 1, $false, "97-1=", 2) | Out-Null
$d.Content.Find.Execute("0+93=", $true, $false, $false, $false, $false, $true, 1, $false, "24+22=", 2) | Out-Null
$d.Content.Find.Execute("42+43=", $true, $false, $false, $false, $false, $true, 1, $false, "11+87=", 2) | Out-Null
$d.Content.Find.Execute("26+33=", $true, $false, $false, $false, $false, $true, 1, $false, "90-10=", 2) | Out-Null
$d.Content.Find.Execute("68-11=", $true, $false, $false, $false, $false, $true, 1, $false, "92-25=", 2) | Out-Null
$d.Content.Find.Execute("51-45=", $true, $false, $false, $false, $false, $true, 1, $false, "41+23=", 2) | Out-Null
$d.Content.Find.Execute("8+49=", $true, $false, $false, $false, $false, $true, 1, $false, "35-4=", 2) | Out-Null
$d.Content.Find.Execute("43-31=", $true, $false, $false, $false, $false, $true, 1, $false, "41+8=", 2) | Out-Null
$d.Content.Find.Execute("73-50=", $true, $false, $false, $false, $false, $true, 1, $false, "74-68=", 2) | Out-Null
$d.Content.Find.Execute("65+1=", $true, $false, $false, $false, $false, $true, 1, $false, "8+43=", 2) | Out-Null
$d.Content.Find.Execute("76+0=", $true, $false, $false, $false, $false, $true, 1, $false, "97-94=", 2) | Out-Null
$d.Content.Find.Execute("78+9=", $true, $false, $false, $false, $false, $true, 1, $false, "67-16=", 2) | Out-Null
$d.Content.Find.Execute("35+64=", $true, $false, $false, $false, $false, $true, 1, $false, "42-6=", 2) | Out-Null
$d.Content.Find.Execute("16+1=", $true, $false, $false, $false, $false, $true, 1, $false, "65-12=", 2) | Out-Null
$d.Content.Find.Execute("0+43=", $true, $false, $false, $false, $false, $true, 1, $false, "92-87=", 2) | Out-Null
$d.Content.Find.Execute("61-17=", $true, $false, $false, $false, $false, $true, 1, $false, "15-9=", 2) | Out-Null
$d.Content.Find.Execute("97-37=", $true, $false, $false, $false, $false, $true, 1, $false, "42+24=", 2) | Out-Null
$d.Content.Find.Execute("2+92=", $true, $false, $false, $false, $false, $true, 1, $false, "34+57=", 2) | Out-Null
$d.Content.Find.Execute("93-41=", $true, $false, $false, $false, $false, $true, 1, $false, "39-14=", 2) | Out-Null
$d.Content.Find.Execute("87-49=", $true, $false, $false, $false, $false, $true, 1, $false, "89-45=", 2) | Out-Null
$d.Content.Find.Execute("40+22=", $true, $false, $false, $false, $false, $true, 1, $false, "55-29=", 2) | Out-Null
$d.Content.Find.Execute("78-19=", $true, $false, $false, $false, $false, $true, 1, $false, "78-48=", 2) | Out-Null
$d.Content.Find.Execute("50+48=", $true, $false, $false, $false, $false, $true, 1, $false, "37+32=", 2) | Out-Null
$d.Content.Find.Execute("45-41=", $true, $false, $false, $false, $false, $true, 1, $false, "23+64=", 2) | Out-Null
$d.Content.Find.Execute("23+22=", $true, $false, $false, $false, $false, $true, 1, $false, "5+56=", 2) | Out-Null
$d.Content.Find.Execute("9+87=", $true, $false, $false, $false, $false, $true, 1, $false, "6+63=", 2) | Out-Null
$d.Content.Find.Execute("80-66=", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=", 2) | Out-Null
$d.Content.Find.Execute("28+35=", $true, $false, $false, $false, $false, $true, 1, $false, "24+14=", 2) | Out-Null
$d.Content.Find.Execute("20+66=", $true, $false, $false, $false, $false, $true, 1, $false, "0+39=", 2) | Out-Null
$d.Content.Find.Execute("91-8=", $true, $false, $false, $false, $false, $true, 1, $false, "35+53=", 2) | Out-Null
$d.Content.Find.Execute("99-8=", $true, $false, $false, $false, $false, $true, 1, $false, "52+31=", 2) | Out-Null
$d.Content.Find.Execute("35+21=", $true, $false, $false, $false, $false, $true, 1, $false, "61-46=", 2) | Out-Null
$d.Content.Find.Execute("11+30=", $true, $false, $false, $false, $false, $true, 1, $false, "15+21=", 2) | Out-Null
$d.Content.Find.Execute("16+7=", $true, $false, $false, $false, $false, $true, 1, $false, "52+15=", 2) | Out-Null
$d.Content.Find.Execute("44-24=", $true, $false, $false, $false, $false, $true, 1, $false, "18-0=", 2) | Out-Null
$d.Content.Find.Execute("65+2=", $true, $false, $false, $false, $false, $true, 1, $false, "1+92=", 2) | Out-Null
$d.Content.Find.Execute("30+61=", $true, $false, $false, $false, $false, $true, 1, $false, "56+43=", 2) | Out-Null
$d.Content.Find.Execute("30+50=", $true, $false, $false, $false, $false, $true, 1, $false, "43-1=", 2) | Out-Null
$d.Content.Find.Execute("46-40=", $true, $false, $false, $false, $false, $true, 1, $false, "29-12=", 2) | Out-Null
$d.Content.Find.Execute("36+2=", $true, $false, $false, $false, $false, $true, 1, $false, "95-7=", 2) | Out-Null
$d.Content.Find.Execute("59-40=", $true, $false, $false, $false, $false, $true, 1, $false, "83+9=", 2) | Out-Null
$d.Content.Find.Execute("98-37=", $true, $false, $false, $false, $false, $true, 1, $false, "28+25=", 2) | Out-Null
$d.Content.Find.Execute("51-42=", $true, $false, $false, $false, $false, $true, 1, $false, "0+3=", 2) | Out-Null
$d.Content.Find.Execute("88+8=", $true, $false, $false, $false, $false, $true, 1, $false, "25+38=", 2) | Out-Null
$d.Content.Find.Execute("47+19=", $true, $false, $false, $false, $false, $true, 1, $false, "68+29=", 2) | Out-Null
$d.Content.Find.Execute("97-63=", $true, $false, $false, $false, $false, $true, 1, $false, "56-33=", 2) | Out-Null
$d.Content.Find.Execute("37+26=", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=", 2) | Out-Null
$d.Content.Find.Execute("35+26=", $true, $false, $false, $false, $false, $true, 1, $false, "67-55=", 2) | Out-Null
$d.Content.Find.Execute("77-32=", $true, $false, $false, $false, $false, $true, 1, $false, "73-28=", 2) | Out-Null
$d.Content.Find.Execute("62-17=", $true, $false, $false, $false, $false, $true, 1, $false, "92-27=", 2) | Out-Null
$d.Content.Find.Execute("89-9=", $true, $false, $false, $false, $false, $true, 1, $false, "35+3=", 2) | Out-Null
$d.Content.Find.Execute("60-11=", $true, $false, $false, $false, $false, $true, 1, $false, "61+14=", 2) | Out-Null
$d.Content.Find.Execute("85+11=", $true, $false, $false, $false, $false, $true, 1, $false, "20-1=", 2) | Out-Null
$d.Content.Find.Execute("13+83=", $true, $false, $false, $false, $false, $true, 1, $false, "29-11=", 2) | Out-Null
$d.Content.Find.Execute("63-35=", $true, $false, $false, $false, $false, $true, 1, $false, "28-15=", 2) | Out-Null
$d.Content.Find.Execute("31-24=", $true, $false, $false, $false, $false, $true, 1, $false, "16+58=", 2) | Out-Null
$d.Content.Find.Execute("3+95=", $true, $false, $false, $false, $false, $true, 1, $false, "70-13=", 2) | Out-Null
$d.Content.Find.Execute("95-73=", $true, $false, $false, $false, $false, $true, 1, $false, "88-48=", 2) | Out-Null
$d.Content.Find.Execute("72+25=", $true, $false, $false, $false, $false, $true, 1, $false, "10+81=", 2) | Out-Null
$d.Content.Find.Execute("43-25=", $true, $false, $false, $false, $false, $true, 1, $false, "59+21=", 2) | Out-Null
$d.Content.Find.Execute("72+20=", $true, $false, $false, $false, $false, $true, 1, $false, "91-76=", 2) | Out-Null
$d.Content.Find.Execute("83-42=", $true, $false, $false, $false, $false, $true, 1, $false, "4+83=", 2) | Out-Null
$d.Content.Find.Execute("38+17=", $true, $false, $false, $false, $false, $true, 1, $false, "57-13=", 2) | Out-Null
$d.Content.Find.Execute("70-57=", $true, $false, $false, $false, $false, $true, 1, $false, "24+73=", 2) | Out-Null
$d.Content.Find.Execute("85-15=", $true, $false, $false, $false, $false, $true, 1, $false, "35+14=", 2) | Out-Null
$d.Content.Find.Execute("76-37=", $true, $false, $false, $false, $false, $true, 1, $false, "28-14=", 2) | Out-Null
$d.Content.Find.Execute("47-7=", $true, $false, $false, $false, $false, $true, 1, $false, "39+59=", 2) | Out-Null
$d.Content.Find.Execute("44-38=", $true, $false, $false, $false, $false, $true, 1, $false, "55+42=", 2) | Out-Null
$d.Content.Find.Execute("18+80=", $true, $false, $false, $false, $false, $true, 1, $false, "41+34=", 2) | Out-Null
$d.Content.Find.Execute("18+78=", $true, $false, $false, $false, $false, $true, 1, $false, "71-69=", 2) | Out-Null
$d.Content.Find.Execute("57-10=", $true, $false, $false, $false, $false, $true, 1, $false, "93-8=", 2) | Out-Null
$d.Content.Find.Execute("92-90=", $true, $false, $false, $false, $false, $true, 1, $false, "78+8=", 2) | Out-Null
$d.Content.Find.Execute("73+7=", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=", 2) | Out-Null
$d.Content.Find.Execute("93-42=", $true, $false, $false, $false, $false, $true, 1, $false, "4+32=", 2) | Out-Null
$d.Content.Find.Execute("73+9=", $true, $false, $false, $false, $false, $true, 1, $false, "14+77=", 2) | Out-Null
$d.Content.Find.Execute("67+14=", $true, $false, $false, $false, $false, $true, 1, $false, "45+1=", 2) | Out-Null
$d.Content.Find.Execute("82-12=", $true, $false, $false, $false, $false, $true, 1, $false, "61-27=", 2) | Out-Null
$d.Content.Find.Execute("12+72=", $true, $false, $false, $false, $false, $true, 1, $false, "13-11=", 2) | Out-Null
$d.Content.Find.Execute("39+51=", $true, $false, $false, $false, $false, $true, 1, $false, "93-9=", 2) | Out-Null
$d.Content.Find.Execute("10+56=", $true, $false, $false, $false, $false, $true, 1, $false, "90-54=", 2) | Out-Null
$d.Content.Find.Execute("26-1=", $true, $false, $false, $false, $false, $true, 1, $false, "48+16=", 2) | Out-Null
